# Updates cryptos list price (column D) and 1h-volume-change (column E) figures
# to match the latest scrape, per the commit "Updated cryptos list on Fri Oct 20
# 21:18:04 UTC 2023 with GitHub Actions".
#
# Numeric-looking price strings in column D are written with a leading
# apostrophe so Excel stores them as text (matching the workbook's existing
# text-typed Price column) instead of silently coercing them to numbers
# (which would also drop meaningful trailing zeros, e.g. "26.80" -> 26.8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.592.63"
$ws.Range("E2").Value = "  +3.25%  "
$ws.Range("D3").Value = "1.604.59"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'212.37"
$ws.Range("D6").Value = "'0.523"
$ws.Range("E6").Value = "  +2.57%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'26.80"
$ws.Range("E8").Value = "  +7.74%  "
$ws.Range("D9").Value = "'43.47"
$ws.Range("E9").Value = "  -4.98%  "
$ws.Range("E10").Value = "  +2.09%  "
$ws.Range("E11").Value = "  +2.57%  "
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "1.833.35"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "1.547.91"
$ws.Range("E14").Value = "  -1.18%  "
$ws.Range("D15").Value = "29.595.77"
$ws.Range("E15").Value = "  +3.16%  "
$ws.Range("D16").Value = "'0.536"
$ws.Range("E16").Value = "  +3.71%  "
$ws.Range("E17").Value = "  +2.23%  "
$ws.Range("E18").Value = "  +3.18%  "
$ws.Range("D19").Value = "'240.68"
$ws.Range("E19").Value = "  +5.30%  "
$ws.Range("E20").Value = "  +3.60%  "
$ws.Range("D21").Value = "0.0₃0692"
$ws.Range("E21").Value = "  +1.67%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "'9.22"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").Value = "'154.27"
$ws.Range("E26").Value = "  +1.64%  "
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("D28").Value = "'15.27"
$ws.Range("E28").Value = "  +3.27%  "
$ws.Range("E29").Value = "  +2.59%  "
$ws.Range("E30").Value = "  -0.18%  "
$ws.Range("E31").Value = "  +3.12%  "
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").Value = "'3.10"
$ws.Range("E34").Value = "  +4.05%  "
$ws.Range("D35").Value = "1.408.35"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E36").Value = "  +0.73%  "
$ws.Range("E37").Value = "  +5.25%  "
$ws.Range("D38").Value = "'2.82"
$ws.Range("E38").Value = "  +3.83%  "
$ws.Range("D39").Value = "'2.30"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("E40").Value = "  +2.62%  "
$ws.Range("D41").Value = "'0.537"
$ws.Range("E41").Value = "  +3.77%  "
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("E43").Value = "  +7.94%  "
$ws.Range("D44").Value = "'53.90"
$ws.Range("E44").Value = "  +27.45%  "
$ws.Range("D45").Value = "'0.797"
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("D46").Value = "'0.998"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "'65.84"
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "1.743.86"
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").Value = "'0.858"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "'86.53"
$ws.Range("E51").Value = "  +2.12%  "
